$d = $word.ActiveDocument

$newText = "Dates à utiliser pour la Campagne Hercules: 13-22 juin, du 12 au 21 juillet, du 10 au 19 août"

# Variant with a trailing space left over from a following (now-removed) run.
$oldWithTrailingSpace = "Dates à utiliser pour la Campagne 2018 Persée:  Du 30 octobre au 8 novembre et du 29 novembre au 8 décembre "
# Variant without the trailing space.
$oldPlain = "Dates à utiliser pour la Campagne 2018 Persée:  Du 30 octobre au 8 novembre et du 29 novembre au 8 décembre"

$rng = $d.Content
$rng.Find.Execute($oldWithTrailingSpace, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
while ($rng.Find.Found) {
    $rng.Delete()
    $rng.InsertAfter($newText)
    $rng.Collapse(0)
    $rng.End = $d.Content.End
    $rng.Find.Execute($oldWithTrailingSpace, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
}

$rng2 = $d.Content
$rng2.Find.Execute($oldPlain, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
while ($rng2.Find.Found) {
    $rng2.Delete()
    $rng2.InsertAfter($newText)
    $rng2.Collapse(0)
    $rng2.End = $d.Content.End
    $rng2.Find.Execute($oldPlain, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
}
